$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that were dropped from the missing-data sample:
# original row 26 ("RM 232") and original row 28 ("SC 92"). Deleting row 26
# first shifts row 28 up to row 27, so delete row 27 next.
$ws.Rows("26").Delete()
$ws.Rows("27").Delete()

# Re-roll which cells are "missing" (blank) vs populated across the
# remaining rows, matching the new missingness pattern.
$ws.Range("C2").Value = 14.9
$ws.Range("E3").ClearContents()
$ws.Range("E4").Value = -6.4
$ws.Range("C6").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("C12").Value = 12.5
$ws.Range("C14").ClearContents()
$ws.Range("E15").Value = -8.4
$ws.Range("E18").Value = -8.5
$ws.Range("E19").ClearContents()
$ws.Range("C20").Value = 12.5
$ws.Range("C21").Value = 12.7
$ws.Range("E22").ClearContents()
$ws.Range("C23").ClearContents()
$ws.Range("E23").Value = -7
$ws.Range("C24").ClearContents()
$ws.Range("E25").Value = -7.1
$ws.Range("B26").Value = -20.2
$ws.Range("B27").ClearContents()
$ws.Range("E27").ClearContents()
$ws.Range("B28").ClearContents()
$ws.Range("B29").Value = -19.5
$ws.Range("B30").Value = -19.7
$ws.Range("B31").ClearContents()
$ws.Range("C31").Value = 15.3
$ws.Range("B32").ClearContents()
$ws.Range("C33").Value = 10.4
